# "muitas mudanças de design"
#
# The "Primeiro Nome" / "Sobrenome" columns are merged into a single
# "Nome completo" column, which shifts every column after it one slot to
# the left (CPF, Data de Nascimento, Email, Cargo/Função all move left,
# and the table + its autofilter + the role dropdown validation shrink
# from A1:F2 to A1:E2). The "Data de Nascimento" column also switches
# from the custom DD/MM/YYYY number format to Excel's built-in short
# date format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet range backing the table is about to be resized, so drop
# the ListObject first (it gets rebuilt afterwards) to avoid it lagging
# behind the column shift.
$lo = $ws.ListObjects.Item(1)
$loName = $lo.Name
$lo.Unlist()

# "Primeiro Nome" (A) becomes "Nome completo"; "Sobrenome" (B) disappears
# entirely and everything to its right shifts left by one column.
$ws.Range("A1").Value = "Nome completo"
$ws.Columns("B").Delete()

# Recreate the "Usuarios" table over the new, narrower A1:E2 range with a
# header row, same as before.
$lo2 = $ws.ListObjects.Add(1, $ws.Range("A1:E2"), 0, 1)
$lo2.Name = $loName

# "Data de Nascimento" is now column C - use Excel's built-in short date
# format instead of the old custom "DD/MM/YYYY" one.
$ws.Columns("C").NumberFormat = "mm-dd-yy"

# Keep the same kind of selection state Excel leaves behind (one cell to
# the right of the table, past the last column).
$null = $ws.Range("G2").Select()
